# Updates cryptos.xlsx "Sheet1" to the Dec 28 2022 03:37 UTC symbol-list refresh.
# - Refreshes Price (column D) figures for rows whose coin lineup didn't change.
# - Rows 10-18 and 42-43 show the coin list shifting: each row's Coin/Link/Price/
#   Volume(1h) cells are replaced with the next entry in the refreshed ranking
#   (row 10 used to be "One", now holds what was row 11's "WazirX", etc.),
#   while the leading rank number embedded in column E stays tied to the row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.76"
$ws.Range("D3").Value = "'23.92"
$ws.Range("D4").Value = "'5.305"
$ws.Range("D5").Value = "'0.05775"
$ws.Range("D6").Value = "'6.480"
$ws.Range("D7").Value = "'3.338"
$ws.Range("D8").Value = "'0.8096"
$ws.Range("D9").Value = "'0.8859"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1394"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07344"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Value = "'0.03120"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03062"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09348"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").Value = "'3.853"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001539"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "'0.04734"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "'0.0006039"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("D19").Value = "'0.005876"
$ws.Range("D20").Value = "'0.001274"
$ws.Range("D22").Value = "'0.00008805"
$ws.Range("D23").Value = "'3.601"
$ws.Range("D28").Value = "'0.0002353"
$ws.Range("D40").Value = "'0.03806"
$ws.Range("D41").Value = "'0.006383"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1055"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.002751"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("D44").Value = "'0.008387"
$ws.Range("D45").Value = "'0.00005400"
$ws.Range("D47").Value = "'0.6911"
$ws.Range("D48").Value = "'0.001846"
$ws.Range("D49").Value = "'0.00002103"
$ws.Range("D50").Value = "'0.0002003"
